# "Running a2,a4 and b6" - update Runmode/Results on the "Test Cases" summary
# sheet for TestCase_A2 and TestCase_A4, and sync the "Results" header cell on
# each TestCase_A* detail sheet to reflect the run.

$wb = $excel.ActiveWorkbook

# --- "Test Cases" summary sheet -------------------------------------------
$ws = $wb.Worksheets.Item("Test Cases")

# TestCase_A1 (row 2) is no longer running
$ws.Range("C2").Value = "N"

# TestCase_A2 (row 3) is now running
$ws.Range("C3").Value = "Y"

# TestCase_A4 (row 5) is now running and passed
$ws.Range("C5").Value = "Y"
$ws.Range("D5").Value = "PASS"

# Leave the selection on C5, matching the last-touched cell
$ws.Range("C5").Select()

# --- detail sheets: update the "Results" header cell -----------------------
$wb.Worksheets.Item("TestCase_A5").Range("D1").Value = "PASS"
$wb.Worksheets.Item("TestCase_A6").Range("D1").Value = "PASS"
$wb.Worksheets.Item("TestCase_A7").Range("D1").Value = "PASS"
$wb.Worksheets.Item("TestCase_A8").Range("D1").Value = "PASS"
$wb.Worksheets.Item("TestCase_A9").Range("D1").Value = "PASS"
$wb.Worksheets.Item("TestCase_A10").Range("F1").Value = "PASS"
$wb.Worksheets.Item("TestCase_A11").Range("D1").Value = "PASS"
$wb.Worksheets.Item("TestCase_A12").Range("F1").Value = "PASS"
